$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos list refresh (prices/volumes) + two row swaps (Fetch.AI/PancakeSwap, Maker/FirstDigitalUSD)

$ws.Range("D2").Value = "64.495.11"
$ws.Range("E2").Value = "  +0.95%  "
$ws.Range("D3").Value = "2.760.60"
$ws.Range("E3").Value = "  +0.31%  "
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "577.22"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.37%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "160.31"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.52%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.603"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.44%  "
$ws.Range("E9").Value = "  -1.31%  "
$ws.Range("E10").Value = "  +4.90%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.81"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.26%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.387"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.73%  "
$ws.Range("D13").Value = "3.251.42"
$ws.Range("E13").Value = "  +0.36%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.28"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.05%  "
$ws.Range("D15").Value = "64.133.29"
$ws.Range("E15").Value = "  +0.46%  "
$ws.Range("E16").Value = "  -1.69%  "
$ws.Range("D17").Value = "2.769.85"
$ws.Range("E17").Value = "  +0.57%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.16"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.88%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.85"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.28%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "358.27"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.11%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.68"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.18%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.998"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.05%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.529"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -7.07%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.16"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.85%  "
$ws.Range("E25").Value = "  -1.09%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.60"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.88%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.19%  "
$ws.Range("D28").Value = "0.0₃0925"
$ws.Range("E28").Value = "  -1.64%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.36"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.36%  "
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.98"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.46%  "
$ws.Range("B31").Value = "Fetch.AI"
$ws.Range("C31").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.38"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +8.97%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "167.66"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.42%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.00"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.06%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.52"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.11%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "20.21"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.84%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.999"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.01%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.84"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.14%  "
$ws.Range("E38").Value = "  -1.22%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "351.65"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.81%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.41"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.00%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.18"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.94%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "39.00"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.34%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "22.52"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.18%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "21.58"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.35%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0592"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.22%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "136.69"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.27%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.631"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.64%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0253"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.01%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.101"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.07%  "
$ws.Range("B50").Value = "FirstDigitalUSD"
$ws.Range("C50").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.997"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.25%  "
$ws.Range("B51").Value = "Maker"
$ws.Range("C51").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D51").Value = "2.146.12"
$ws.Range("E51").Value = "  +0.85%  "
